# Auto-generated edit script for THREAT_ALERT sheet update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Step 1: New row data (rows 2..23 after the edit). Captured first so
# we can write it after formats are settled.
# -----------------------------------------------------------------
$rowData = @(
    @{ Row=2; A="05-MAR-26"; B="SM-486"; C="EgyptAir MS-680"; D=713; E=560; F=153; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=3; A="08-MAR-26"; B="SM-486"; C="EgyptAir MS-680"; D=713; E=560; F=153; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=4; A="12-MAR-26"; B="SM-486"; C="EgyptAir MS-680"; D=948; E=758; F=190; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=5; A="15-MAR-26"; B="SM-486"; C="flyadeal F3-775"; D=797; E=1088; F=-291; G=30; H=30; I=0; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=6; A="15-MAR-26"; B="SM-486"; C="EgyptAir MS-676"; D=1003; E=1088; F=-85; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=7; A="15-MAR-26"; B="SM-486"; C="EgyptAir MS-696"; D=1003; E=1088; F=-85; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=8; A="15-MAR-26"; B="SM-486"; C="EgyptAir MS-640"; D=1003; E=1088; F=-85; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=9; A="15-MAR-26"; B="SM-486"; C="EgyptAir MS-680"; D=1078; E=1088; F=-10; G=46; H=30; I=-16; J="MEDIUM THREAT - MONITOR"; JStyle=4; K="SAR" },
    @{ Row=10; A="17-MAR-26"; B="SM-486"; C="EgyptAir MS-694"; D=1278; E=1298; F=-20; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=11; A="17-MAR-26"; B="SM-486"; C="EgyptAir MS-696"; D=1278; E=1298; F=-20; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=12; A="17-MAR-26"; B="SM-486"; C="EgyptAir MS-640"; D=1278; E=1298; F=-20; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=13; A="17-MAR-26"; B="SM-486"; C="EgyptAir MS-678"; D=1335; E=1298; F=37; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=14; A="24-MAY-26"; B="SM-486"; C="Saudia SV-381"; D=679; E=914; F=-235; G=46; H=30; I=-16; J="MEDIUM THREAT - MONITOR"; JStyle=4; K="SAR" },
    @{ Row=15; A="24-MAY-26"; B="SM-486"; C="Saudia SV-319"; D=679; E=914; F=-235; G=46; H=30; I=-16; J="MEDIUM THREAT - MONITOR"; JStyle=4; K="SAR" },
    @{ Row=16; A="24-MAY-26"; B="SM-486"; C="Saudia SV-391"; D=679; E=914; F=-235; G=46; H=30; I=-16; J="MEDIUM THREAT - MONITOR"; JStyle=4; K="SAR" },
    @{ Row=17; A="24-MAY-26"; B="SM-486"; C="flyadeal F3-771"; D=799; E=914; F=-115; G=40; H=30; I=-10; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=18; A="24-MAY-26"; B="SM-486"; C="EgyptAir MS-680"; D=1078; E=914; F=164; G=46; H=30; I=-16; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=19; A="28-MAY-26"; B="SM-486"; C="flyadeal F3-771"; D=619; E=1298; F=-679; G=15; H=30; I=15; J="MEDIUM THREAT - MONITOR"; JStyle=4; K="SAR" },
    @{ Row=20; A="28-MAY-26"; B="SM-486"; C="Saudia SV-391"; D=679; E=1298; F=-619; G=46; H=30; I=-16; J="HIGH THREAT ALERT - NEED ACTION"; JStyle=5; K="SAR" },
    @{ Row=21; A="28-MAY-26"; B="SM-486"; C="Saudia SV-381"; D=679; E=1298; F=-619; G=46; H=30; I=-16; J="HIGH THREAT ALERT - NEED ACTION"; JStyle=5; K="SAR" },
    @{ Row=22; A="28-MAY-26"; B="SM-486"; C="flynas XY-793"; D=889; E=1298; F=-409; G=20; H=30; I=10; J="LOW THREAT"; JStyle=3; K="SAR" },
    @{ Row=23; A="28-MAY-26"; B="SM-486"; C="flynas XY-576"; D=1079; E=1298; F=-219; G=40; H=30; I=-10; J="LOW THREAT"; JStyle=3; K="SAR" }
)

# -----------------------------------------------------------------
# Step 2: propagate the J-column (IMPACT) cell formatting for each
# threat level using pristine donor cells that already carry the
# matching style in the original sheet, BEFORE any values/styles
# are overwritten:
#   style 3 = "LOW THREAT"                   -> donor J2
#   style 4 = "MEDIUM THREAT - MONITOR"       -> donor J6
#   style 5 = "HIGH THREAT ALERT - NEED ACTION" -> donor J24
# -----------------------------------------------------------------
$donors = @{ 3 = "J2"; 4 = "J6"; 5 = "J24" }

foreach ($styleId in $donors.Keys) {
    $donorAddr = $donors[$styleId]
    $ws.Range($donorAddr).Copy()
    foreach ($item in $rowData) {
        if ($item.JStyle -eq $styleId) {
            $ws.Range("J" + $item.Row).PasteSpecial(-4122) | Out-Null
        }
    }
}
$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# Step 3: write the new cell values
# -----------------------------------------------------------------
foreach ($item in $rowData) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
    $ws.Range("H$r").Value = $item.H
    $ws.Range("I$r").Value = $item.I
    $ws.Range("J$r").Value = $item.J
    $ws.Range("K$r").Value = $item.K
}

# -----------------------------------------------------------------
# Step 4: the refreshed report only has 22 data rows (2..23); the
# old trailing rows 24:25 are no longer part of the report.
# -----------------------------------------------------------------
$ws.Rows("24:25").Delete()

